$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-121 down to 43-122
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with data
$ws.Cells.Item(42, 1).Value = 2
$ws.Cells.Item(42, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = 45028
$ws.Cells.Item(42, 5).Value = 4
$ws.Cells.Item(42, 6).Value = 100112030
$ws.Cells.Item(42, 7).Value = "Poroto granado"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 700
$ws.Cells.Item(42, 11).Value = 30000
$ws.Cells.Item(42, 12).Value = 31000
$ws.Cells.Item(42, 13).Value = 30500
$ws.Cells.Item(42, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 1220
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
